# Generate Report for Handoff
# Replace the old e2e test-doc GUID/filename (b5eb0be3-...) with the new one
# (8d142fa8-...) everywhere it shows up in the localization-status report,
# update the handoff-generation timestamps, and refresh the generated
# xliff file names (new content hash) for both locales.

$wb = $excel.ActiveWorkbook

$oldGuidMd   = "b5eb0be3-66eb-4eb9-af13-a87922dd377b.md"
$newGuidMd   = "8d142fa8-10d9-420a-8baa-1aa05f6bca52.md"
$oldPathMd   = "e2e\b5eb0be3-66eb-4eb9-af13-a87922dd377b.md"
$newPathMd   = "e2e\8d142fa8-10d9-420a-8baa-1aa05f6bca52.md"

$oldHoDate   = "2016-08-12 23:12:47"
$newHoDate   = "2016-08-12 23:13:21"

$oldZhHandoffFile = "b5eb0be3-66eb-4eb9-af13-a87922dd377b.494a858dad02ff056956000b2130605900c3f316.zh-cn.xlf"
$newZhHandoffFile = "8d142fa8-10d9-420a-8baa-1aa05f6bca52.72edb2ab0556393848e9ad449f04c7741384edba.zh-cn.xlf"
$oldZhHandoffDate = "2016-08-12 23:12:39"
$newZhHandoffDate = "2016-08-12 23:13:14"

$oldDeHandoffFile = "b5eb0be3-66eb-4eb9-af13-a87922dd377b.494a858dad02ff056956000b2130605900c3f316.de-de.xlf"
$newDeHandoffFile = "8d142fa8-10d9-420a-8baa-1aa05f6bca52.72edb2ab0556393848e9ad449f04c7741384edba.de-de.xlf"

# Hyperlink targets (relationship URLs) stay pointed at the old commit/file;
# only the cell text + hyperlink display text change.
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/oltest/blob/a647c710635fc925f61206bffcc9099df2bb60bc/e2e/b5eb0be3-66eb-4eb9-af13-a87922dd377b.md"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuidMd
$wsOverview.Range("B2").Value = $newPathMd
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", $newPathMd)

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newGuidMd
$wsZhCn.Range("G2").Value = $newZhHandoffFile
$wsZhCn.Range("H2").Value = $newZhHandoffDate

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, "", "", $newGuidMd)

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newGuidMd
$wsDeDe.Range("G2").Value = $newDeHandoffFile
$wsDeDe.Range("H2").Value = $newHoDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, "", "", $newGuidMd)
